# Refactor timetable name accurately
# Corrects Semester 2/3 module rows and replaces the stray "OSS1014" entry
# with the real module code "CSC3044 " (trailing space preserved, as in source).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Semester 2 block (rows 6-9): SEG1201, NET1014, WEB1201, PRG1203
$ws.Range("B7").Value = "NET1014"
$ws.Range("B8").Value = "WEB1201"
$ws.Range("B9").Value = "PRG1203"

# Semester 3 block (rows 10-13): WEB1201, SEG1201, NET1014, CSC3044
$ws.Range("B10").Value = "WEB1201"
$ws.Range("B11").Value = "SEG1201"
$ws.Range("B12").Value = "NET1014"
$ws.Range("B13").Value = "CSC3044 "

# Update the active selection to match the saved view state
$ws.Range("R32").Select()
